$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value-only updates (column E) ---
$ws.Range("E7").Value  = 11
$ws.Range("E13").Value = 23
$ws.Range("E14").Value = 11
$ws.Range("E15").Value = 10
$ws.Range("E19").Value = 15
$ws.Range("E22").Value = 12
$ws.Range("E23").Value = 11
$ws.Range("E24").Value = 13
$ws.Range("E28").Value = 17

# --- Row 29: was "La morte di Silvio Berlusconi"/FanPage/Instagram/16 ---
#             now "L'incoronazione di Re Carlo"/La Repubblica/YouTube/3
$ws.Range("B29").Value = "L'incoronazione di Re Carlo"
$ws.Range("C29").Value = "La Repubblica"
$ws.Range("D29").Value = "YouTube"
$ws.Range("E29").Value = 3

# --- Row 30: D Youtube -> Facebook, E 5 -> 1 ---
$ws.Range("D30").Value = "Facebook"
$ws.Range("E30").Value = 1

# --- Row 31: C Il Corriere Della Sera -> FanPage, D Facebook -> Instagram, E 3 -> 17 ---
$ws.Range("C31").Value = "FanPage"
$ws.Range("D31").Value = "Instagram"
$ws.Range("E31").Value = 17

# --- Row 32: C Il Corriere Della Sera -> FanPage, D Instagram -> YouTube, E 24 -> 5 ---
$ws.Range("C32").Value = "FanPage"
$ws.Range("D32").Value = "YouTube"
$ws.Range("E32").Value = 5

# --- Row 33: D YouTube -> Facebook (E stays 4) ---
$ws.Range("D33").Value = "Facebook"

# --- Row 34: C La Repubblica -> Il Corriere Della Sera, E 17 -> 25 ---
$ws.Range("C34").Value = "Il Corriere Della Sera"
$ws.Range("E34").Value = 25

# --- Row 35: C La Repubblica -> Il Corriere Della Sera (D stays YouTube), E 8 -> 4 ---
$ws.Range("C35").Value = "Il Corriere Della Sera"
$ws.Range("E35").Value = 4

# --- New row 36: "La morte di Silvio Berlusconi"/La Repubblica/Instagram/17 ---
$ws.Range("A36").Value = ""
$ws.Range("B36").Value = "La morte di Silvio Berlusconi"
$ws.Range("C36").Value = "La Repubblica"
$ws.Range("D36").Value = "Instagram"
$ws.Range("E36").Value = 17

# --- New row 37: "La morte di Silvio Berlusconi"/La Repubblica/YouTube/12 ---
$ws.Range("A37").Value = ""
$ws.Range("B37").Value = "La morte di Silvio Berlusconi"
$ws.Range("C37").Value = "La Repubblica"
$ws.Range("D37").Value = "YouTube"
$ws.Range("E37").Value = 12
